# Underline the phrase "this form" inside the sentence:
#   "...please fill out this form (not available in demo mode) so that..."
#
# Before:  <w:r><w:t xml:space="preserve">t this form </w:t></w:r>
# After:   <w:r><w:t xml:space="preserve">t </w:t></w:r>
#          <w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>this form</w:t></w:r>
#          <w:r><w:t xml:space="preserve"> </w:t></w:r>
#
# Word automatically splits the owning run into the run(s) needed to carry
# the new direct formatting when only part of a run's text is reformatted,
# which is exactly the shape seen in the diff.
#
# (The diff's other hunk only touches wp14:editId on an unrelated picture's
# <wp:anchor> - a random, content-free bookkeeping stamp Word mints on save
# that isn't exposed anywhere in the Word object model/COM surface, so it
# isn't something any automation script can target deliberately.)

$d = $word.ActiveDocument

# "this form" is unique in the document body, so a simple Find on the
# document's main Content range locates the correct phrase unambiguously.
$target = $d.Content
$found = $target.Find.Execute("this form", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the phrase 'this form' to underline."
}

# $target now spans exactly "this form" (Find.Execute collapses/extends the
# range it was called on to the match). Apply single underline formatting to
# just that sub-range; Word splits the surrounding run(s) as needed so the
# new formatting applies only to this text, leaving the rest of the
# paragraph's runs/text untouched.
$target.Font.Underline = 1
